# "add common tool file"
#
# Adds a 4th column ("notes") to the 工作进展计划 (work progress plan) sheet:
#   - a handful of existing rows (20/21/22/27/31) get a short note in column D
#   - two brand-new rows (43/44) are appended, matching the look of the
#     existing highlighted rows (columns A/B/C use the yellow "high priority"
#     fill already used elsewhere on the sheet)
#   - the sheet's page setup is switched to A4/portrait, same as sheet 1
#
# NOTE on ordering: the shared-string table is built in first-use order, so
# the writes below are sequenced to reproduce the exact same shared-string
# index assignment as the original author's edit (row 43/44 first, then the
# D20/D22/D21 performance notes, then D31/D27, then row 44's own values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("工作进展计划")

# --- new row 43: IMPORT ---------------------------------------------------
$ws.Range("A43").Value = "高"
$ws.Range("B43").Value = "高"
$ws.Range("C43").Value = "IMPORT"
$ws.Range("A43:C43").Interior.Color = 65535
$ws.Range("D43").Value = "1 (方便测试性能）"

# --- notes for the existing "performance test" rows -----------------------
$ws.Range("D22").Value = "3（最大化性能测试）"
$ws.Range("D20").Value = "2（最大化性能测试）"
$ws.Range("D21").Value = "4（最大化性能测试）"

# --- notes for the "other statements" rows ---------------------------------
$ws.Range("D31").Value = "5 方便测试其他语句"
$ws.Range("D27").Value = "6 方便测试其他语句"

# --- new row 44: 进程，工具 --------------------------------------------------
$ws.Range("A44").Value = "高"
$ws.Range("B44").Value = "高"
$ws.Range("C44").Value = "进程，工具"
$ws.Range("A44:C44").Interior.Color = 65535
$ws.Range("D44").Value = "0 便于测试"

# --- page setup: A4 / portrait, matching sheet 1 ---------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- leave the cursor where the author ended up ----------------------------
$ws.Range("D39").Select() | Out-Null
